# JGT.xlsx - "Lecture partielle de l'EDT M1 MIAGE."
#
# The timetable's day-name labels are renamed (jeudi/vendredi/samedi ->
# lundi/mardi/mercredi) and the schedule dates are shifted forward by
# 1096 days (3 years) to the new matching weekday.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the weekday labels used throughout the sheet (shared text, so a
# single replace updates every occurrence).
$ws.Cells.Replace("jeudi", "lundi")
$ws.Cells.Replace("vendredi", "mardi")
$ws.Cells.Replace("samedi", "mercredi")

# Shift the date serials in column A forward by 1096 days (3 years) to
# align with the new weekday labels.
$dateMap = @{
  "A2"  = 46027
  "A4"  = 46028
  "A7"  = 46029
  "A10" = 46034
  "A12" = 46035
  "A15" = 46042
  "A18" = 46050
  "A21" = 46057
  "A24" = 46092
  "A27" = 46104
  "A30" = 46111
  "A32" = 46113
  "A35" = 46147
  "A38" = 46153
}

foreach ($addr in $dateMap.Keys) {
    $ws.Range($addr).Value2 = $dateMap[$addr]
}
